# Change the "E3-value" labels to "e3value" with the "3" rendered as a
# superscript (matching the e3value notation), on both AutoShape 47
# shapes that carry this label (inside the two "Grouper 69" groups).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Set-E3ValueLabel($shape) {
    $tr = $shape.TextFrame.TextRange
    $paraCount = $tr.Paragraphs().Count
    for ($pIdx = 1; $pIdx -le $paraCount; $pIdx++) {
        $para = $tr.Paragraphs($pIdx, 1)
        # Non-last paragraphs carry a trailing paragraph-mark (CR) in
        # their .Text, so trim before comparing.
        if ($para.Text.Trim() -eq "E3-value") {
            $para.Text = "e3value"
            # Make the "3" (2nd character) superscript; PowerPoint
            # automatically splits the run into e / 3 / value pieces.
            $three = $para.Characters(2, 1)
            $three.Font.Superscript = $true
        }
    }
}

$shapeCount = $s.Shapes.Count
for ($shapeIdx = 1; $shapeIdx -le $shapeCount; $shapeIdx++) {
    $top = $s.Shapes.Item($shapeIdx)
    if ($top.Name -eq "Grouper 69") {
        $itemCount = $top.GroupItems.Count
        for ($itemIdx = 1; $itemIdx -le $itemCount; $itemIdx++) {
            $item = $top.GroupItems.Item($itemIdx)
            if ($item.HasTextFrame) {
                Set-E3ValueLabel $item
            }
        }
    }
}
